# Apply the edits described by the commit: update the wording of the
# "purpose" / "intent" blurbs on the Instructions sheet, resize their
# rows, and update the saved selections on both sheets.

$wb = $excel.ActiveWorkbook

$wsInstructions = $wb.Worksheets.Item("Instructions")
$wsParameters   = $wb.Worksheets.Item("parameters")

# --- Instructions sheet: updated text -------------------------------------

$purposeText = 'The purpose of this workbook is to define a group of ports. Although the intent of port groups was to examine port usage of server clusters or storage enclosures, any groups of ports can be defined. The examples on the "parameters" sheet was used for test. It was left in as an example of how to do everything. Typically, SAN admins use a simple alias naming convention so that it''s easy to identify a server or storage cluster by and alias prefix. In that case, you would only name the group,  Set Filter to "Alias", the Operand to "some_prefix_*", and the Operator to "Wild".'

$intentText = 'This workbook is intended for use with the -group option in report.py and stats_g.py.'

$wsInstructions.Range("A1").Value = $purposeText
$wsInstructions.Range("A3").Value = $intentText

# Row height changes that accompany the new wording
$wsInstructions.Range("A1:C1").RowHeight = 90
$wsInstructions.Range("A3:C3").RowHeight = 15

# --- Selection / view state -------------------------------------------------

$wsParameters.Activate()
$wsParameters.Range("B10").Select()

$wsInstructions.Activate()
$wsInstructions.Range("A1:C1").Select()
